$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (Spanish labels -> short English codes) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case fixes for connector words (de/del/la/las/el/los/y -> De/Del/La/Las/El/Los/Y) ---
$ws.Range("B5").Value = "Pabellón De Arteaga"
$ws.Range("B6").Value = "Rincón De Romos"
$ws.Range("B7").Value = "San José De Gracia"
$ws.Range("B23").Value = "Amatenango De La Frontera"
$ws.Range("B59").Value = "Guadalupe Y Calvo"
$ws.Range("B61").Value = "Hidalgo Del Parral"
$ws.Range("B76").Value = "Valle De Zaragoza"
$ws.Range("B83").Value = "San Juan De Sabinas"
$ws.Range("B91").Value = "Villa De Álvarez"
$ws.Range("A93").Value = "Ciudad De México"
$ws.Range("B110").Value = "Coneto De Comonfort"
$ws.Range("B123").Value = "Nombre De Dios"
$ws.Range("B126").Value = "Pánuco De Coronado"
$ws.Range("B132").Value = "San Juan De Guadalupe"
$ws.Range("A142").Value = "Estado De México"
$ws.Range("B142").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B143").Value = "Almoloya De Alquisiras"
$ws.Range("B144").Value = "Almoloya De Juárez"
$ws.Range("B147").Value = "Atizapán De Zaragoza"
$ws.Range("B152").Value = "Coacalco De Berriozábal"
$ws.Range("B156").Value = "Ecatepec De Morelos"
$ws.Range("B157").Value = "Ixtapan De La Sal"
$ws.Range("B162").Value = "Naucalpan De Juárez"
$ws.Range("B165").Value = "San Antonio La Isla"
$ws.Range("B166").Value = "San Felipe Del Progreso"
$ws.Range("B170").Value = "Tenango Del Valle"
$ws.Range("B172").Value = "Tlalnepantla De Baz"
$ws.Range("B176").Value = "Valle De Bravo"
$ws.Range("B177").Value = "Valle De Chalco Solidaridad"
$ws.Range("B178").Value = "Villa De Allende"
$ws.Range("B187").Value = "Apaseo El Alto"
$ws.Range("B193").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B197").Value = "Jaral Del Progreso"
$ws.Range("B204").Value = "Purísima Del Rincón"
$ws.Range("B208").Value = "San Diego De La Unión"
$ws.Range("B210").Value = "San Francisco Del Rincón"
$ws.Range("B212").Value = "San Luis De La Paz"
$ws.Range("B213").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B217").Value = "Valle De Santiago"
$ws.Range("B220").Value = "Acapulco De Juárez"
$ws.Range("B222").Value = "Alcozauca De Guerrero"
$ws.Range("B225").Value = "Atoyac De Álvarez"
$ws.Range("B226").Value = "Ayutla De Los Libres"
$ws.Range("B229").Value = "Chilpancingo De Los Bravo"
$ws.Range("B230").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B233").Value = "Coyuca De Benítez"
$ws.Range("B237").Value = "Cuetzala Del Progreso"
$ws.Range("B238").Value = "Cutzamala De Pinzón"
$ws.Range("B244").Value = "Huitzuco De Los Figueroa"
$ws.Range("B245").Value = "Iguala De La Independencia"
$ws.Range("B247").Value = "Zihuatanejo De Azueta"
$ws.Range("B259").Value = "Taxco De Alarcón"
$ws.Range("B261").Value = "Técpan De Galeana"
$ws.Range("B263").Value = "Tixtla De Guerrero"
$ws.Range("B266").Value = "Tlapa De Comonfort"
$ws.Range("B272").Value = "Atotonilco El Grande"
$ws.Range("B274").Value = "Cuautepec De Hinojosa"
$ws.Range("B276").Value = "Huasca De Ocampo"
$ws.Range("B280").Value = "Mixquiahuala De Juárez"
$ws.Range("B281").Value = "Pachuca De Soto"
$ws.Range("B283").Value = "Progreso De Obregón"
$ws.Range("B286").Value = "Tepehuacán De Guerrero"
$ws.Range("B289").Value = "Tula De Allende"
$ws.Range("B290").Value = "Tulancingo De Bravo"
$ws.Range("B292").Value = "Zacualtipán De Ángeles"
$ws.Range("B297").Value = "Atemajac De Brizuela"
$ws.Range("B299").Value = "Atotonilco El Alto"
$ws.Range("B300").Value = "Autlán De Navarro"
$ws.Range("B314").Value = "Huejuquilla El Alto"
$ws.Range("B321").Value = "La Manzanilla De La Paz"
$ws.Range("B322").Value = "Lagos De Moreno"
$ws.Range("B326").Value = "Ojuelos De Jalisco"
$ws.Range("B328").Value = "San Cristóbal De La Barranca"
$ws.Range("B329").Value = "San Diego De Alejandría"
$ws.Range("B330").Value = "San Juan De Los Lagos"
$ws.Range("B333").Value = "San Miguel El Alto"
$ws.Range("B336").Value = "Tamazula De Gordiano"
$ws.Range("B341").Value = "Teocuitatlán De Corona"
$ws.Range("B342").Value = "Tepatitlán De Morelos"
$ws.Range("B343").Value = "Tizapán El Alto"
$ws.Range("B349").Value = "Unión De San Antonio"
$ws.Range("B352").Value = "Zapotlán Del Rey"
$ws.Range("B353").Value = "Zapotlán El Grande"
$ws.Range("B422").Value = "Puente De Ixtla"
$ws.Range("B424").Value = "Tlaltizapán De Zapata"
$ws.Range("B430").Value = "Bahía De Banderas"
$ws.Range("B450").Value = "Chalcatongo De Hidalgo"
$ws.Range("B452").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B455").Value = "Mariscala De Juárez"
$ws.Range("B457").Value = "Nejapa De Madero"
$ws.Range("B458").Value = "Oaxaca De Juárez"
$ws.Range("B459").Value = "Putla Villa De Guerrero"
$ws.Range("B461").Value = "San Agustín De Las Juntas"
$ws.Range("B463").Value = "San Antonino El Alto"
$ws.Range("B465").Value = "San Dionisio Del Mar"
$ws.Range("B468").Value = "San Juan Del Estado"
$ws.Range("B478").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B480").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B500").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B501").Value = "Tlacolula De Matamoros"
$ws.Range("B502").Value = "Villa De Etla"
$ws.Range("B503").Value = "Zimatlán De Álvarez"
$ws.Range("B519").Value = "Cuayuca De Andrade"
$ws.Range("B523").Value = "Izúcar De Matamoros"
$ws.Range("B525").Value = "Los Reyes De Juárez"
$ws.Range("B526").Value = "Palmar De Bravo"
$ws.Range("B534").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B537").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B548").Value = "Amealco De Bonfil"
$ws.Range("B549").Value = "Cadereyta De Montes"
$ws.Range("B552").Value = "Jalpan De Serra"
$ws.Range("B554").Value = "Pinal De Amoles"
$ws.Range("B556").Value = "San Juan Del Río"
$ws.Range("B565").Value = "Ciudad Del Maíz"
$ws.Range("B572").Value = "Santa María Del Río"
$ws.Range("B577").Value = "Villa De Ramos"
$ws.Range("B626").Value = "Soto La Marina"
$ws.Range("B635").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B641").Value = "Amatlán De Los Reyes"
$ws.Range("B649").Value = "Cosamaloapan De Carpio"
$ws.Range("B653").Value = "Ixhuatlán De Madero"
$ws.Range("B660").Value = "Martínez De La Torre"
$ws.Range("B690").Value = "Jiménez Del Teul"
$ws.Range("B693").Value = "Nochistlán De Mejía"
$ws.Range("B700").Value = "Villa De Cos"

# --- Floating point last-digit correction (1 ULP) for 3/3195 percentage cells ---
$ws.Range("D13").Value = 0.0009389671361502348
$ws.Range("D14").Value = 0.0009389671361502348
$ws.Range("D40").Value = 0.0009389671361502348
$ws.Range("D41").Value = 0.0009389671361502348
$ws.Range("D60").Value = 0.0009389671361502348
$ws.Range("D67").Value = 0.0009389671361502348
$ws.Range("D71").Value = 0.0009389671361502348
$ws.Range("D81").Value = 0.0009389671361502348
$ws.Range("D83").Value = 0.0009389671361502348
$ws.Range("D104").Value = 0.0009389671361502348
$ws.Range("D115").Value = 0.0009389671361502348
$ws.Range("D127").Value = 0.0009389671361502348
$ws.Range("D140").Value = 0.0009389671361502348
$ws.Range("D158").Value = 0.0009389671361502348
$ws.Range("D161").Value = 0.0009389671361502348
$ws.Range("D174").Value = 0.0009389671361502348
$ws.Range("D180").Value = 0.0009389671361502348
$ws.Range("D192").Value = 0.0009389671361502348
$ws.Range("D193").Value = 0.0009389671361502348
$ws.Range("D209").Value = 0.0009389671361502348
$ws.Range("D243").Value = 0.0009389671361502348
$ws.Range("D244").Value = 0.0009389671361502348
$ws.Range("D248").Value = 0.0009389671361502348
$ws.Range("D252").Value = 0.0009389671361502348
$ws.Range("D254").Value = 0.0009389671361502348
$ws.Range("D255").Value = 0.0009389671361502348
$ws.Range("D258").Value = 0.0009389671361502348
$ws.Range("D278").Value = 0.0009389671361502348
$ws.Range("D281").Value = 0.0009389671361502348
$ws.Range("D295").Value = 0.0009389671361502348
$ws.Range("D306").Value = 0.0009389671361502348
$ws.Range("D315").Value = 0.0009389671361502348
$ws.Range("D332").Value = 0.0009389671361502348
$ws.Range("D339").Value = 0.0009389671361502348
$ws.Range("D340").Value = 0.0009389671361502348
$ws.Range("D342").Value = 0.0009389671361502348
$ws.Range("D345").Value = 0.0009389671361502348
$ws.Range("D347").Value = 0.0009389671361502348
$ws.Range("D349").Value = 0.0009389671361502348
$ws.Range("D353").Value = 0.0009389671361502348
$ws.Range("D357").Value = 0.0009389671361502348
$ws.Range("D368").Value = 0.0009389671361502348
$ws.Range("D371").Value = 0.0009389671361502348
$ws.Range("D384").Value = 0.0009389671361502348
$ws.Range("D393").Value = 0.0009389671361502348
$ws.Range("D400").Value = 0.0009389671361502348
$ws.Range("D444").Value = 0.0009389671361502348
$ws.Range("D467").Value = 0.0009389671361502348
$ws.Range("D475").Value = 0.0009389671361502348
$ws.Range("D488").Value = 0.0009389671361502348
$ws.Range("D501").Value = 0.0009389671361502348
$ws.Range("D505").Value = 0.0009389671361502348
$ws.Range("D532").Value = 0.0009389671361502348
$ws.Range("D560").Value = 0.0009389671361502348
$ws.Range("D587").Value = 0.0009389671361502348
$ws.Range("D601").Value = 0.0009389671361502348
$ws.Range("D602").Value = 0.0009389671361502348
$ws.Range("D610").Value = 0.0009389671361502348
$ws.Range("D611").Value = 0.0009389671361502348
$ws.Range("D613").Value = 0.0009389671361502348
$ws.Range("D616").Value = 0.0009389671361502348
$ws.Range("D623").Value = 0.0009389671361502348
$ws.Range("D627").Value = 0.0009389671361502348
$ws.Range("D650").Value = 0.0009389671361502348
$ws.Range("D682").Value = 0.0009389671361502348
$ws.Range("D696").Value = 0.0009389671361502348
$ws.Range("D698").Value = 0.0009389671361502348
$ws.Range("D700").Value = 0.0009389671361502348

# --- Remove trailing footer/metadata rows (708:712) ---
$ws.Rows("708:712").Delete()
